$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert three new rows at 30-32 (existing rows 30-32 shift down to 33-35)
# ---------------------------------------------------------------------------
$ws.Rows("30:32").Insert()

# ---------------------------------------------------------------------------
# 2. Copy formatting onto the freshly inserted rows from rows that already
#    carry (close to) the desired look, then patch up the few cells whose
#    formatting differs from their donor row.
# ---------------------------------------------------------------------------
$ws.Range("A22:I22").Copy()
$ws.Range("A30:I30").PasteSpecial(-4122)

$ws.Range("A23:I23").Copy()
$ws.Range("A31:I31").PasteSpecial(-4122)

$ws.Range("A13:I13").Copy()
$ws.Range("A32:I32").PasteSpecial(-4122)

# Row 32 has a couple of cells whose shading differs from row 13's:
$ws.Range("B10").Copy()
$ws.Range("B32").PasteSpecial(-4122)

$ws.Range("C13").Copy()
$ws.Range("D32").PasteSpecial(-4122)

# The pre-existing last row (now row 35) picks up the border used by the row
# above it once the new rows have been inserted above it.
$ws.Range("D34").Copy()
$ws.Range("D35").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Fill in the values for the three new rows.
# ---------------------------------------------------------------------------
$d45044 = Get-Date -Year 2023 -Month 4 -Day 28 -Hour 0 -Minute 0 -Second 0

$ws.Range("A30").Value = $d45044
$ws.Range("B30").Value = "RASY"
$ws.Range("C30").Value = 178
$ws.Range("D30").Value = "Lotbinière"
$ws.Range("E30").Value = "Chaudière-Appalaches"
$ws.Range("F30").Value = "B"
$ws.Range("G30").Value = "Cote 3"
$ws.Range("H30").Value = "Summum d'activité"
$ws.Range("I30").Value = "Bertrand Le Grand"

$ws.Range("A31").Value = $d45044
$ws.Range("B31").Value = "PSCR"
$ws.Range("C31").Value = 178
$ws.Range("D31").Value = "Lotbinière"
$ws.Range("E31").Value = "Chaudière-Appalaches"
$ws.Range("F31").Value = "B"
$ws.Range("G31").Value = "Cote 3"
$ws.Range("H31").Value = "Summum d'activité"
$ws.Range("I31").Value = "Bertrand Le Grand"

$ws.Range("A32").Value = $d45044
$ws.Range("B32").Value = "BUAM"
$ws.Range("C32").Value = "95m"
$ws.Range("D32").Value = "Saint-Joachim"
$ws.Range("E32").Value = "Capitale-Nationale"
$ws.Range("F32").Value = "C/D"
$ws.Range("G32").Value = "Cote 1-2"
$ws.Range("H32").Value = "Réserve nationale de faune du Cap-Tourmente"
$ws.Range("I32").Value = "Simon Bourbeau"

# ---------------------------------------------------------------------------
# 4. Restore the sort range to cover the now-larger table and move the
#    active selection, matching the saved view state.
# ---------------------------------------------------------------------------
[void]$ws.Range("A4:I35").Sort($ws.Range("A4:A35"), 1, $ws.Range("F4:F35"), [System.Reflection.Missing]::Value, 1, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 1)

[void]$ws.Range("M25").Select()
